$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (changed) date column C for rows 2-10 from 2023-09-01 (45170)
# to 2023-09-05 (45174), keeping the existing date serial / formatting.
$ws.Range("C2:C10").Value = 45174
